# "removing future course attempts from prediction"
#
# This adds a new day's log entry (2024-05-15) to the working-hours sheet,
# bumps the hours recorded for the previous day, introduces a new
# "TODO" column (F) and records a short note about the next task.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix up the previous day's hours (B13: 5 -> 8) ---
$ws.Cells.Item(13, 2).Value2 = 8

# --- new row 14: 2024-05-15 (serial 45427), 4 hours ---
# A14: copy the date style from A13 so it keeps the short-date format
$ws.Cells.Item(13, 1).Copy() | Out-Null
$ws.Cells.Item(14, 1).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Cells.Item(14, 1).Value2 = 45427

# B14: plain hours value, default style like the other Hours cells
$ws.Cells.Item(14, 2).Value2 = 4

# D14: note about removing future course attempts - written first so it
# lands before the new "TODO" string in the shared-strings table
$ws.Cells.Item(14, 4).WrapText = $true
$ws.Cells.Item(14, 4).Value2 = "removed all ""future"" courses (aka retries at the same course AND courses of a higher year level)"

# --- new column F: "TODO" header + note on the new row ---
$ws.Cells.Item(1, 5).Copy() | Out-Null
$ws.Cells.Item(1, 6).PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Cells.Item(1, 6).Value2 = "TODO"

$ws.Cells.Item(14, 6).WrapText = $true
$ws.Cells.Item(14, 6).Value2 = "compare relevancies between courses| create a multi output regression model"

# row height for the new (wrapped) row, matching row 12's wrapped style
$ws.Rows.Item(14).RowHeight = 42.75

# give the new column a sensible width
$ws.Columns.Item(6).ColumnWidth = 30.42

# page orientation, as recorded in the saved workbook
$ws.PageSetup.Orientation = 1

# leave the selection where the author left it after the edit
$null = $ws.Range("F15").Select()
